$d = $word.ActiveDocument
$apos = [char]8217

# 1. Title: " crowdfunding data based " -> " crowdfunding data "
$d.Content.Find.Execute(
    " crowdfunding data based ", $true, $false, $false, $false, $false,
    $true, 1, $false, " crowdfunding data ", 2)

# 2. First red question: rewrite whole sentence.
$d.Content.Find.Execute(
    "Given the provided data, what are three conclusions we can draw about Kickstarter campaigns?",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Three conclusions could be drawn from Kickstarter campaigns?", 2)

# 3. "theater campaigns" bullet: We -> I, analyze -> analyzed, have -> had
$d.Content.Find.Execute(
    "When we analyze the data using a Stacked Column Pivot Chart based on category, theater campaigns have the most success.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "When I analyzed the data using a Stacked Column Pivot Chart based on category, theater campaigns had the most success.",
    2)

# 4. "plays" bullet: We -> I, have -> had
$d.Content.Find.Execute(
    "When we analyze the data using a Stacked Column Pivot Chart based on sub-category, play",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "When I analyze the data using a Stacked Column Pivot Chart based on sub-category, play",
    2)

$d.Content.Find.Execute(
    "s have the most success within t",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "s had the most success within t",
    2)

# 5. Second red question: rewrite whole sentence.
$d.Content.Find.Execute(
    "What are some limitations of this dataset?",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Some limitations of this data set?", 2)

# 6. "We don't know why projects are more successful " -> "I don't know why projects are more successful "
$d.Content.Find.Execute(
    "We don" + $apos + "t know why projects are more successful ",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "I don" + $apos + "t know why projects are more successful ", 2)

# 7. "We don't know wh" (followed by "at makes...") -> "I don't know wh"
$d.Content.Find.Execute(
    "We don" + $apos + "t know what makes a project successful",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "I don" + $apos + "t know what makes a project successful", 2)

# 8. "We don't know why " (followed by "certain countries...") -> "I don't know why "
$d.Content.Find.Execute(
    "We don" + $apos + "t know why certain countries favor specific campaigns",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "I don" + $apos + "t know why certain countries favor specific campaigns", 2)
